$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 231, shifting existing rows 231-246 down to 232-247
$ws.Rows("231:231").Insert()

# Populate the new row 231 with the new weekly data point
$ws.Range("A231").Value = 4
$ws.Range("B231").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C231").Value = "Los Lagos"
$ws.Range("D231").Value = 44585
$ws.Range("E231").Value = 10
$ws.Range("F231").Value = 100112045
$ws.Range("G231").Value = "Zapallo"
$ws.Range("H231").Value = "Paine"
$ws.Range("I231").Value = "1a nueva(o)"
$ws.Range("J231").Value = 500
$ws.Range("K231").Value = 500
$ws.Range("L231").Value = 500
$ws.Range("M231").Value = 500
$ws.Range("N231").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O231").Value = "Región de O'Higgins"
$ws.Range("P231").Value = 500
$ws.Range("Q231").Value = 1
$ws.Range("R231").Value = "Hortaliza"
